# Apply updated crypto price/volume figures to Sheet1 (generated from commit diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.904.38'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '2.338.14'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''306.60'
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("D6").Value = '''100.34'
$ws.Range("E6").Value = '  -1.44%  '
$ws.Range("D7").Value = '''0.512'
$ws.Range("E7").Value = '  -4.56%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -3.53%  '
$ws.Range("D10").Value = '''34.93'
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("D11").Value = '''52.12'
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("D14").Value = '''6.79'
$ws.Range("E14").Value = '  -2.80%  '
$ws.Range("D15").Value = '''15.95'
$ws.Range("E15").Value = '  +6.69%  '
$ws.Range("D16").Value = '2.286.07'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = '''0.805'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '42.834.78'
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").Value = '''6.21'
$ws.Range("E19").Value = '  +0.82%  '
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("D21").Value = '''11.74'
$ws.Range("E21").Value = '  -4.55%  '
$ws.Range("D22").Value = '''67.96'
$ws.Range("D23").Value = '''236.88'
$ws.Range("E23").Value = '  -1.83%  '
$ws.Range("D24").Value = '''2.02'
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").Value = '''2.56'
$ws.Range("E25").Value = '  -2.29%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = '''25.55'
$ws.Range("E27").Value = '  +3.83%  '
$ws.Range("D28").Value = '''3.95'
$ws.Range("E28").Value = '  -1.05%  '
$ws.Range("E29").Value = '  +9.61%  '
$ws.Range("D30").Value = '''35.02'
$ws.Range("E30").Value = '  -4.66%  '
$ws.Range("D31").Value = '''9.32'
$ws.Range("E31").Value = '  -3.32%  '
$ws.Range("D32").Value = '''160.22'
$ws.Range("E32").Value = '  -4.65%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").Value = '''5.12'
$ws.Range("E34").Value = '  -2.88%  '
$ws.Range("D35").Value = '''4.67'
$ws.Range("E35").Value = '  +7.44%  '
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").Value = '''17.38'
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("D38").Value = '''0.0728'
$ws.Range("E39").Value = '  -3.71%  '
$ws.Range("D40").Value = '''1.86'
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("E41").Value = '  -3.11%  '
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("D43").Value = '''2.43'
$ws.Range("E43").Value = '  +5.05%  '
$ws.Range("D44").Value = '2.014.98'
$ws.Range("E44").Value = '  +2.36%  '
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("D46").Value = '''18.74'
$ws.Range("E46").Value = '  -2.47%  '
$ws.Range("D47").Value = '''10.29'
$ws.Range("E47").Value = '  +3.75%  '
$ws.Range("E48").Value = '  -1.80%  '
$ws.Range("D49").Value = '''55.89'
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").Value = '2.563.50'
$ws.Range("E51").Value = '  +1.09%  '
